$wb = $excel.ActiveWorkbook

# Duplicate the "Grading" sheet to use as the starting point for the new
# "2023" sheet (inherits all formatting: column widths, row heights, styles).
$grading = $wb.Worksheets.Item("Grading")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$grading.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2023"

# Remove the "Runs on Azure without crashing" row - it no longer appears in
# the 2023 rubric. Deleting the row shifts everything below it up by one.
$newSheet.Rows.Item(9).Delete()

# Update the scores that changed for the 2023 revision.
$newSheet.Range("B5").Value = 10

$newSheet.Range("B6").Value = 15
$newSheet.Range("D6").Value = 15
$newSheet.Range("E6").Value = 15

$newSheet.Range("D9").Value = 7
$newSheet.Range("E9").Value = 7

$newSheet.Range("D13").Value = 8
$newSheet.Range("E13").Value = 8

# Selection on the new sheet.
$newSheet.Range("A2:E16").Select()
